# Update the "Förändrad" (Changed) date column (C) for rows 2-120
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C120")
$range.Value = 45174
